{"js": "// Append the Bibliograf\u00eda section: two blank paragraphs, a \"Bibliograf\u00eda:\"\n// heading paragraph, and a paragraph with the reference URL \u2014 all inserted\n// after the last (empty) paragraph of the document, inheriting that\n// paragraph's formatting (majorHAnsi theme font, sz/szCs 24).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\n// Two empty paragraphs.\nanchor = anchor.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nanchor = anchor.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// \"Bibliograf\u00eda:\" paragraph.\nanchor = anchor.insertParagraph(\"Bibliograf\u00eda:\", \"After\");\nawait context.sync();\n\n// Paragraph with the reference link.\nanchor = anchor.insertParagraph(\n  \"https://www.hektorprofe.net/tutorial/django-formularios-crear-editar-instancias\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "# Append the Bibliograf\u00eda section: two blank paragraphs, a \"Bibliograf\u00eda:\"\n# heading paragraph, and a paragraph with the reference URL \u2014 all inserted\n# after the last (empty) paragraph of the document, inheriting that\n# paragraph's formatting (majorHAnsi theme font, sz/szCs 24).\n$d = $word.ActiveDocument\n\n# Two empty paragraphs.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n\n# \"Bibliograf\u00eda:\" paragraph.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"Bibliograf\u00eda:\"\n\n# Paragraph with the reference link.\n$d.Paragraphs.Last.Range.InsertParagraphAfter()\n$d.Paragraphs.Last.Range.Text = \"https://www.hektorprofe.net/tutorial/django-formularios-crear-editar-instancias\"\n"}
